# Apply the commit "Ajout liens classes aad360e5fee5029c63105e72d3ce4f0f9c39d95f":
#  - Update the "Date" metadata value on the Metadata sheet.
#  - Add a new element row (SituationExercice.exerciceProfessionnel) on the
#    Elements sheet, describing a Reference link to the ExerciceProfessionnel class.

$wb = $excel.ActiveWorkbook

$wsMeta = $wb.Worksheets.Item("Metadata")
$wsElem = $wb.Worksheets.Item("Elements")

# --- Metadata sheet: bump the generation Date value (row 8, column B) ---
$wsMeta.Range("B8").Value = "2025-07-24T13:17:05+00:00"

# --- Elements sheet: append a new row describing the new link element ---
# Start from a duplicate of the last existing row (18): this brings over the
# formatting (style, borders, wrap text) *and* every cell value, including
# the columns that stay identical between row 18 and row 19 (Slice Name,
# Alias(s), Label, Must Support?, Is Modifier?, Is Summary?, Comments,
# Requirements, Default Value, Meaning When Missing, Fixed Value ... up to
# Constraint(s)) so only the cells that actually change need to be touched.
$wsElem.Range("A18:AJ18").Copy($wsElem.Range("A19:AJ19"))

$wsElem.Range("A19").Value = "SituationExercice.exerciceProfessionnel"
$wsElem.Range("B19").Value = "SituationExercice.exerciceProfessionnel"

# Min / Max / Base Min / Base Max are textual "1" values (like the existing
# "0"/"*" cells); force text formatting first so they are not reinterpreted
# as numbers.
$wsElem.Range("F19").NumberFormat = "@"
$wsElem.Range("G19").NumberFormat = "@"
$wsElem.Range("AG19").NumberFormat = "@"
$wsElem.Range("AH19").NumberFormat = "@"
$wsElem.Range("F19").Value = "1"
$wsElem.Range("G19").Value = "1"

$wsElem.Range("K19").Value = "Reference(https://interop.esante.gouv.fr/ig/fhir/mos/StructureDefinition/ExerciceProfessionnel)`n"
$wsElem.Range("L19").Value = "Lien vers la classe ExerciceProfessionnel."
$wsElem.Range("M19").Value = "Lien vers la classe ExerciceProfessionnel."
$wsElem.Range("AF19").Value = "SituationExercice.exerciceProfessionnel"
$wsElem.Range("AG19").Value = "1"
$wsElem.Range("AH19").Value = "1"

# Re-fit the row height: assigning values touched the wrap-text autofit
# height; put it back to the sheet default (same as every other row).
$wsElem.Rows.Item(19).AutoFit()

# Column K (Type(s)) now needs to fit the longer "Reference(...)" text.
$wsElem.Columns.Item(11).ColumnWidth = 74.666666667
